# Scheduled runner update: refresh Leve profit calculations (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69 (Leve Item ID 12616)
$ws.Range("H69").Value = 3690.75
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3690.75
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 11072.25
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -12820.25
# Row 72 (Leve Item ID 12616)
$ws.Range("H72").Value = 3690.75
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3690.75
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 33216.75
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -41952.75
# Row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 2601
$ws.Range("I92").Value = 1124
$ws.Range("J92").Value = 5555
$ws.Range("K92").Value = 1124
$ws.Range("L92").Value = 5555
$ws.Range("M92").Value = 124
$ws.Range("N92").Value = -8051
# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 66679292
$ws.Range("I96").Value = 5119.875
$ws.Range("J96").Value = 142878350
$ws.Range("K96").Value = 15359.625
$ws.Range("L96").Value = 428635050
$ws.Range("M96").Value = -13986.625
$ws.Range("N96").Value = -428637796
# Row 133 (Leve Item ID 41856)
$ws.Range("H133").Value = 35530.77
$ws.Range("J133").Value = 35530.77
$ws.Range("L133").Value = 35530.77
$ws.Range("N133").Value = -45650.77

$ws = $wb.Worksheets.Item("ARM")
# Row 56 (Leve Item ID 2504)
$ws.Range("H56").Value = 12000
$ws.Range("J56").Value = 12000
$ws.Range("L56").Value = 12000
$ws.Range("N56").Value = -13484
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 1368.0667
$ws.Range("I97").Value = 1250.7142
$ws.Range("K97").Value = 1250.7142
$ws.Range("M97").Value = -754.7141999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 75 (Leve Item ID 11872)
$ws.Range("H75").Value = 15623.077
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 15623.077
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 15623.077
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -17495.077
# Row 78 (Leve Item ID 11872)
$ws.Range("H78").Value = 15623.077
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 15623.077
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 46869.231
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -56229.231
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2467.0144
$ws.Range("I134").Value = 1570.6154
$ws.Range("J134").Value = 3594.742
$ws.Range("K134").Value = 4711.8462
$ws.Range("L134").Value = 10784.226
$ws.Range("M134").Value = -2176.8462
$ws.Range("N134").Value = -15854.226

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4553.7793
$ws.Range("I31").Value = 1976.2084
$ws.Range("J31").Value = 5720.981
$ws.Range("K31").Value = 1976.2084
$ws.Range("L31").Value = 5720.981
$ws.Range("M31").Value = -1681.2084
$ws.Range("N31").Value = -6310.981
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4553.7793
$ws.Range("I34").Value = 1976.2084
$ws.Range("J34").Value = 5720.981
$ws.Range("K34").Value = 1976.2084
$ws.Range("L34").Value = 5720.981
$ws.Range("M34").Value = -1774.2084
$ws.Range("N34").Value = -6124.981
# Row 81 (Leve Item ID 10613)
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# Row 84 (Leve Item ID 10613)
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 69 (Leve Item ID 12850)
$ws.Range("H69").Value = 85786600
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 85786600
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 257359800
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -257361422
# Row 72 (Leve Item ID 12850)
$ws.Range("H72").Value = 85786600
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 85786600
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 772079400
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -772087512
# Row 133 (Leve Item ID 44073)
$ws.Range("H133").Value = 5076
$ws.Range("I133").Value = 3988.182
$ws.Range("J133").Value = 6571.75
$ws.Range("K133").Value = 11964.546
$ws.Range("L133").Value = 19715.25
$ws.Range("M133").Value = -6904.545999999998
$ws.Range("N133").Value = -29835.25

$ws = $wb.Worksheets.Item("GSM")
# Row 34 (Leve Item ID 10924)
$ws.Range("H34").Value = 20000
$ws.Range("J34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("N34").Value = -20536
# Row 63 (Leve Item ID 11048)
$ws.Range("H63").Value = 10000
$ws.Range("J63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("N63").Value = -11372
# Row 66 (Leve Item ID 11048)
$ws.Range("H66").Value = 10000
$ws.Range("J66").Value = 10000
$ws.Range("L66").Value = 30000
$ws.Range("N66").Value = -36864
# Row 76 (Leve Item ID 10924)
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630
# Row 79 (Leve Item ID 10924)
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 5040.857
$ws.Range("I80").Value = 5254.2856
$ws.Range("J80").Value = 4827.4287
$ws.Range("K80").Value = 5254.2856
$ws.Range("L80").Value = 4827.4287
$ws.Range("M80").Value = -4256.2856
$ws.Range("N80").Value = -6823.4287
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 5040.857
$ws.Range("I83").Value = 5254.2856
$ws.Range("J83").Value = 4827.4287
$ws.Range("K83").Value = 26271.428
$ws.Range("L83").Value = 24137.1435
$ws.Range("M83").Value = -21279.428
$ws.Range("N83").Value = -34121.14350000001

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 3662.7273
$ws.Range("I40").Value = 2721.7856
$ws.Range("J40").Value = 5309.375
$ws.Range("K40").Value = 2721.7856
$ws.Range("L40").Value = 5309.375
$ws.Range("M40").Value = -2585.7856
$ws.Range("N40").Value = -5581.375
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2855.0557
$ws.Range("I132").Value = 2065.0435
$ws.Range("J132").Value = 4252.769
$ws.Range("K132").Value = 6195.130500000001
$ws.Range("L132").Value = 12758.307
$ws.Range("M132").Value = -3665.130500000001
$ws.Range("N132").Value = -17818.307
# Row 133 (Leve Item ID 41903)
$ws.Range("H133").Value = 36191.668
$ws.Range("J133").Value = 36191.668
$ws.Range("L133").Value = 36191.668
$ws.Range("N133").Value = -41251.668

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1176658.2
$ws.Range("I132").Value = 1450375.2
$ws.Range("J132").Value = 3585.1428
$ws.Range("K132").Value = 4351125.6
$ws.Range("L132").Value = 10755.4284
$ws.Range("M132").Value = -4348595.6
$ws.Range("N132").Value = -15815.4284
